# JIRAGITHUB-7 additional small changes based on feedback from Mar 25
#
# 1) Bump the "last updated" date placeholder (3/22/2020 -> 3/26/2020)
#    on the slide master and every slide layout.
# 2) Slide 2 title: "What is "Git Integration for Jira"?" ->
#    "What is "Git Integration"?"
# 3) Slide 4 body: append ", IMO" to the GitHub/Atlassian footnote.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders ("3/22/2020" -> "3/26/2020") on master + layouts
# ---------------------------------------------------------------------
$newDate = "3/26/2020"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 2 title text change
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "What is*") {
        $shp.TextFrame.TextRange.Text = "What is " + [char]0x201C + "Git Integration" + [char]0x201D + "?"
    }
}

# ---------------------------------------------------------------------
# 3) Slide 4 body paragraph: append ", IMO"
# ---------------------------------------------------------------------
# NB: the COM text-range getter normalizes curly quotes/apostrophes to
# their straight ASCII equivalents, so the target text is rebuilt from
# explicit character codes instead of round-tripping through .Text.
$ldq  = [char]0x201C   # “
$rdq  = [char]0x201D   # ”
$rsq  = [char]0x2019   # '
$finalGitHubText = "** GitHub has " + $ldq + "issue" + $rdq + " functionality, but not a strong customizable workflow model, and their integration marketplace is not yet as mature as Atlassian" + $rsq + "s, IMO"

$slide4 = $p.Slides.Item(4)
for ($i = 1; $i -le $slide4.Shapes.Count; $i++) {
    $shp = $slide4.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "*GitHub has*") {
        $tr = $shp.TextFrame.TextRange
        $allParas = $tr.Paragraphs(0, -1)
        for ($j = 1; $j -le $allParas.Count; $j++) {
            $para = $tr.Paragraphs($j, 1)
            if ($para.Text -like "*GitHub has*") {
                # Two-step write: first blow away the old text with an
                # unrelated placeholder so the run-diffing logic doesn't
                # try to keep a shared-prefix run (which would splinter
                # the paragraph into two runs and inject a stray CR).
                # Then write the real final text in one shot so it lands
                # as a single run with the paragraph's original rPr.
                $para.Text = "zzz_temp_zzz"
                $para2 = $tr.Paragraphs($j, 1)
                $para2.Text = $finalGitHubText
                break
            }
        }
    }
}
